$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: NO. OF HOURS LATE (column E) corrected from 0.25 to 0.5
$ws.Range("E10").Value = 0.5

# Row 18 (05-08-2015, OB Others day): fill in the Official Business
# time-out/time-in details that were previously left blank.
$ws.Range("K18").Value = "08:30:00"
$ws.Range("L18").Value = "08:30:00"
$ws.Range("M18").Value = "18:30:00"
$ws.Range("N18").Value = "18:30:00"
